# Update statistic page + reload data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (drop "năm" from the title)
$ws.Name = "Doanh thu 2023"

# Reload the monthly revenue figures (values refreshed + two new months filled in)
$ws.Range("A2").Value = "Tháng 1"
$ws.Range("B2").Value = 245000

$ws.Range("A3").Value = "Tháng 2"
$ws.Range("B3").Value = 300000

$ws.Range("A4").Value = "Tháng 3"
$ws.Range("B4").Value = 395000

$ws.Range("A5").Value = "Tháng 4"
$ws.Range("B5").Value = 390000

$ws.Range("A6").Value = "Tháng 5"
$ws.Range("B6").Value = 455000

$ws.Range("A7").Value = "Tháng 6"
$ws.Range("B7").Value = 345000

$ws.Range("A8").Value = "Tháng 7"
$ws.Range("B8").Value = 430000

$ws.Range("A9").Value = "Tháng 8"
$ws.Range("B9").Value = 260000

$ws.Range("A10").Value = "Tháng 9"
$ws.Range("B10").Value = 485000

$ws.Range("A11").Value = "Tháng 10"
$ws.Range("B11").Value = 175000

$ws.Range("A12").Value = "Tháng 11"
$ws.Range("B12").Value = 245000

$ws.Range("A13").Value = "Tháng 12"
$ws.Range("B13").Value = 400000

# Grand total now lives on row 15, leaving row 14 blank
$ws.Range("A15").Value = "TỔNG CỘNG"
$ws.Range("B15").Value = 4125000
